$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> (date, ebitda) updates per the diff.
# Row 38 only changes the date (column A); EBITDA (column B) is unchanged.
$updates = @(
    @{ Row = 2;  Date = "2026/01/06"; Ebitda = "7.67" },
    @{ Row = 8;  Date = "2026/01/06"; Ebitda = "8.69" },
    @{ Row = 14; Date = "2026/01/06"; Ebitda = "3.14" },
    @{ Row = 20; Date = "2026/01/06"; Ebitda = "13.00" },
    @{ Row = 26; Date = "2026/01/06"; Ebitda = "11.50" },
    @{ Row = 32; Date = "2026/01/06"; Ebitda = "27.76" },
    @{ Row = 38; Date = "2026/01/06"; Ebitda = $null },
    @{ Row = 44; Date = "2026/01/06"; Ebitda = "12.41" },
    @{ Row = 50; Date = "2026/01/06"; Ebitda = "11.44" },
    @{ Row = 56; Date = "2026/01/06"; Ebitda = "30.78" },
    @{ Row = 62; Date = "2026/01/06"; Ebitda = "11.17" },
    @{ Row = 68; Date = "2026/01/06"; Ebitda = "12.87" },
    @{ Row = 74; Date = "2026/01/06"; Ebitda = "17.67" }
)

foreach ($u in $updates) {
    $cellA = $ws.Cells.Item($u.Row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $u.Date
    $cellA.Style = "Normal"

    if ($u.Ebitda -ne $null) {
        $cellB = $ws.Cells.Item($u.Row, 2)
        $cellB.NumberFormat = "@"
        $cellB.Value = $u.Ebitda
        $cellB.Style = "Normal"
    }
}
